$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.01542466666666666
$ws.Range("H2").Value = 0.046274
$ws.Range("I2").Value = 0.001050900616877799
$ws.Range("J2").Value = 0.001050900616877799
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 1.706626889997333
$ws.Range("R2").Value = 15.359642009976
$ws.Range("S2").Value = 0.0005755172026969992
$ws.Range("T2").Value = 0.0005755172026969993
$ws.Range("G3").Value = 0.01542466666666666
$ws.Range("H3").Value = 0.046274
$ws.Range("I3").Value = 0.001050900616877799
$ws.Range("J3").Value = 0.001050900616877799
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 0.9819748211655555
$ws.Range("R3").Value = 8.83777339049
$ws.Range("S3").Value = 0.0003311464301356284
$ws.Range("T3").Value = 0.0003311464301356284
$ws.Range("G4").Value = 0.01542466666666666
$ws.Range("H4").Value = 0.046274
$ws.Range("I4").Value = 0.001050900616877799
$ws.Range("J4").Value = 0.001050900616877799
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 0.4277173894195554
$ws.Range("R4").Value = 3.849456504775999
$ws.Range("S4").Value = 0.0001442369840451713
$ws.Range("T4").Value = 0.0001442369840451713
$ws.Range("I5").Value = 0.8505602471689909
$ws.Range("J5").Value = 0.8505602471689909
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 1381.280937577156
$ws.Range("R5").Value = 12431.5284381944
$ws.Range("S5").Value = 0.4658024234777735
$ws.Range("T5").Value = 0.4658024234777735
$ws.Range("I6").Value = 0.8505602471689909
$ws.Range("J6").Value = 0.8505602471689909
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("S6").Value = 0.2680177220773685
$ws.Range("T6").Value = 0.2680177220773685
$ws.Range("I7").Value = 0.8505602471689909
$ws.Range("J7").Value = 0.8505602471689909
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 346.1786991276227
$ws.Range("R7").Value = 3115.608292148604
$ws.Range("S7").Value = 0.1167401016138489
$ws.Range("T7").Value = 0.1167401016138489
$ws.Range("G8").Value = 2.177987666666667
$ws.Range("H8").Value = 6.533963
$ws.Range("I8").Value = 0.1483888522141314
$ws.Range("J8").Value = 0.1483888522141314
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 240.9784534306013
$ws.Range("R8").Value = 2168.806080875412
$ws.Range("S8").Value = 0.08126395185818588
$ws.Range("T8").Value = 0.08126395185818588
$ws.Range("G9").Value = 2.177987666666667
$ws.Range("H9").Value = 6.533963
$ws.Range("I9").Value = 0.1483888522141314
$ws.Range("J9").Value = 0.1483888522141314
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 138.6564193375839
$ws.Range("R9").Value = 1247.907774038255
$ws.Range("S9").Value = 0.04675840692588238
$ws.Range("T9").Value = 0.04675840692588238
$ws.Range("G10").Value = 2.177987666666667
$ws.Range("H10").Value = 6.533963
$ws.Range("I10").Value = 0.1483888522141314
$ws.Range("J10").Value = 0.1483888522141314
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 60.39438122755688
$ws.Range("R10").Value = 543.5494310480119
$ws.Range("S10").Value = 0.0203664934300631
$ws.Range("T10").Value = 0.0203664934300631
